# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages update at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 270
    3  = 286
    5  = 848
    8  = 8593
    12 = 113
    19 = 726
    20 = 33
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
